# Applies the "scikit-learn Dynamic Time Warping SVM implementation" commit.
#
# The change breaks several runs in two so that a <w:proofErr w:type="gramStart"/>
# / <w:proofErr w:type="gramEnd"/> pair brackets a short "grammar-flagged"
# phrase (e.g. "doesn't", "has to", "take into account", "Similar to"), merges
# three runs of one paragraph back into a single run, and appends three new
# list paragraphs (two with content, one empty) at the end of the document.
#
# Because <w:proofErr/> markers are structural (not visible text), they can't
# be produced through Find/Replace; instead each affected paragraph's whole
# Range is rewritten in one shot via Range.InsertXML with the exact OOXML for
# that paragraph (preserving its original w:pPr).

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Paragraph 7: "Although PCA yielded good results, it doesn't matter ..."
# Split "doesn't" into its own proofErr-wrapped run.
# ---------------------------------------------------------------------------
$p7 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Although PCA yielded good results, it </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>doesn't</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> matter if the columns will be jumbled, we will get the same result. </w:t></w:r>
</w:p>
"@
$d.Paragraphs(7).Range.InsertXML($p7)

# ---------------------------------------------------------------------------
# Paragraph 8: "A bias has to be added ... [As seen in the above plot ...]"
# Split "has to" (first sentence) and "doesn't" (inside the bracketed,
# Helvetica-styled quote) into their own proofErr-wrapped runs.
# ---------------------------------------------------------------------------
$helv = '<w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/>'
$helvShd = $helv + '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>'

$p8 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">A bias </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>has to</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> be added to the columns towards the previous columns as they are dependent (time series).</w:t></w:r>
  <w:r><w:t xml:space="preserve"> [</w:t></w:r>
  <w:r><w:rPr>$helvShd</w:rPr><w:t xml:space="preserve">As seen in the above plot, the data with shuffled columns also yields similar accuracy scores as PCA </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr>$helvShd</w:rPr><w:t>doesn't</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr>$helvShd</w:rPr><w:t xml:space="preserve"> consider the order (order matters here as the data is a time series) which is flawed for this particular dataset.</w:t></w:r>
  <w:r><w:rPr>$helv</w:rPr><w:br/></w:r>
  <w:r><w:rPr>$helv</w:rPr><w:br/></w:r>
  <w:r><w:rPr>$helvShd</w:rPr><w:t xml:space="preserve">So, we need to consider the bias of a column towards </w:t></w:r>
  <w:r><w:rPr>$helvShd</w:rPr><w:t>its</w:t></w:r>
  <w:r><w:rPr>$helvShd</w:rPr><w:t xml:space="preserve"> previous columns.</w:t></w:r>
  <w:r><w:t>]</w:t></w:r>
</w:p>
"@
$d.Paragraphs(8).Range.InsertXML($p8)

# ---------------------------------------------------------------------------
# Paragraph 12: "T" + "he sigmoid ... 1" + ". " + "T" + "he probability ... classes."
# Split "take into account" into its own proofErr-wrapped run.
# ---------------------------------------------------------------------------
$p12 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>
  </w:pPr>
  <w:r><w:t>T</w:t></w:r>
  <w:r><w:t>he sigmoid activation function gives the value between 0 and 1</w:t></w:r>
  <w:r><w:t xml:space="preserve">. </w:t></w:r>
  <w:r><w:t>T</w:t></w:r>
  <w:r><w:t xml:space="preserve">he probability that the data point belongs to class 1 does not </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>take into account</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> the probability of the other classes.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(12).Range.InsertXML($p12)

# ---------------------------------------------------------------------------
# Paragraph 13: "Similar to the sigmoid activation function ..." + picture.
# Split "Similar to" (the very first run) into its own proofErr-wrapped run.
# ---------------------------------------------------------------------------
# NB: kept on one line (no incidental whitespace) - w:drawing's children are
# stored as an opaque blob, so stray indentation would become literal text.
$drawingXml = '<w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="7A07B724" wp14:editId="34FF9226"><wp:extent cx="2197100" cy="616143"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="256519549" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="256519549" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2212509" cy="620464"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing>'

$p13 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr>
    <w:jc w:val="center"/>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>Similar to</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> the sigmoid activation function the SoftMax function returns the probability of each class. Here is the equation for the SoftMax activation function.</w:t></w:r>
  <w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:rPr><w:noProof/></w:rPr>$drawingXml</w:r>
</w:p>
"@
$d.Paragraphs(13).Range.InsertXML($p13)

# ---------------------------------------------------------------------------
# Paragraph 14: "After the prediction, ... standard scalar has to be applied ..."
# Split "has to" into its own proofErr-wrapped run (strike-through formatting).
# ---------------------------------------------------------------------------
$p14 = @"
<w:p $wns>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">After the prediction, the inverse transform of the standard scalar </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>has to</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> be applied (which was used to reduce bias) to get the actual </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:rPr><w:strike/></w:rPr><w:t>y_predicted</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(14).Range.InsertXML($p14)

# ---------------------------------------------------------------------------
# Paragraph 23 (last paragraph): merge its three Helvetica runs into one run,
# then append three new list paragraphs after it (two with new "My ideas" /
# "Used Scikit-learn ..." content, one empty) in the same InsertXML call so
# the trailing w:sectPr stays put.
# ---------------------------------------------------------------------------
$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

$p23 = @"
<w:p $wns>
  $listPPr
  <w:r><w:rPr>$helvShd</w:rPr><w:t>The simulation of causality or creating the bias towards the previous data is working as expected as shuffling of the columns did not yield the same result as the unshuffled data. The accuracy is very low (94% compared to 5.3% in this case) meaning that a dependency has been created between the columns.</w:t></w:r>
</w:p>
<w:p $wns>
  $listPPr
  <w:r><w:t>My ideas: Give the series as a vector [Check by iterating through various sizes of vector i.e., divide the whole vector into 2, 3</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t xml:space="preserve"> ,...</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> n(number of columns)]</w:t></w:r>
  <w:r><w:t xml:space="preserve">. Did not work as arrays are given for training and so the </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>doesn&#8217;t</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> matter the size of the vector, everything is converted to arrays and each is broken into individual features (which is the data itself!).</w:t></w:r>
</w:p>
<w:p $wns>
  $listPPr
  <w:r><w:t xml:space="preserve">Used Scikit-learn </w:t></w:r>
  <w:r><w:t>Dynamic Time Warping SVM</w:t></w:r>
  <w:r><w:t xml:space="preserve"> and it </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>doesn&#8217;t</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> stop running.</w:t></w:r>
</w:p>
<w:p $wns>
  $listPPr
</w:p>
"@
$d.Paragraphs(23).Range.InsertXML($p23)
